$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 13334097
$ws.Range("I19").Value = 41666990
$ws.Range("J19").Value = 970.94116
$ws.Range("K19").Value = 41666990
$ws.Range("L19").Value = 970.94116
$ws.Range("M19").Value = -41666815
$ws.Range("N19").Value = -1320.94116

$ws.Range("H32").Value = 573.2857
$ws.Range("I32").Value = 567.75
$ws.Range("J32").Value = 580.6667
$ws.Range("K32").Value = 567.75
$ws.Range("L32").Value = 580.6667
$ws.Range("M32").Value = -241.75
$ws.Range("N32").Value = -1232.6667

$ws.Range("H33").Value = 205.1875
$ws.Range("I33").Value = 184.24
$ws.Range("J33").Value = 280
$ws.Range("K33").Value = 184.24
$ws.Range("L33").Value = 280
$ws.Range("M33").Value = 44.75999999999999
$ws.Range("N33").Value = -738

$ws.Range("H63").Value = 40271
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 40271
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41519

$ws.Range("H64").Value = 368173.56
$ws.Range("I64").Value = 570347.8
$ws.Range("J64").Value = 4259.9
$ws.Range("K64").Value = 570347.8
$ws.Range("L64").Value = 4259.9
$ws.Range("M64").Value = -570099.8
$ws.Range("N64").Value = -4755.9

$ws.Range("H66").Value = 40271
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 40271
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127053

$ws.Range("H67").Value = 368173.56
$ws.Range("I67").Value = 570347.8
$ws.Range("J67").Value = 4259.9
$ws.Range("K67").Value = 570347.8
$ws.Range("L67").Value = 4259.9
$ws.Range("M67").Value = -569489.8
$ws.Range("N67").Value = -5975.9

$ws.Range("H80").Value = 1181.5
$ws.Range("I80").Value = 606.6667
$ws.Range("J80").Value = 2139.5557
$ws.Range("K80").Value = 1820.0001
$ws.Range("L80").Value = 6418.6671
$ws.Range("M80").Value = -822.0001
$ws.Range("N80").Value = -8414.667099999999

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H83").Value = 1181.5
$ws.Range("I83").Value = 606.6667
$ws.Range("J83").Value = 2139.5557
$ws.Range("K83").Value = 5460.0003
$ws.Range("L83").Value = 19256.0013
$ws.Range("M83").Value = -468.0002999999997
$ws.Range("N83").Value = -29240.0013

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H96").Value = 396.2857
$ws.Range("I96").Value = 396.2857
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1188.8571
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 184.1428999999998
$ws.Range("N96").ClearContents()

$ws.Range("H137").Value = 2464.6428
$ws.Range("I137").Value = 1697.7142
$ws.Range("J137").Value = 6299.2856
$ws.Range("K137").Value = 5093.142599999999
$ws.Range("L137").Value = 18897.8568
$ws.Range("M137").Value = -2543.142599999999
$ws.Range("N137").Value = -23997.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1147.2759
$ws.Range("I2").Value = 1165.6818
$ws.Range("J2").Value = 1089.4286
$ws.Range("K2").Value = 1165.6818
$ws.Range("L2").Value = 1089.4286
$ws.Range("M2").Value = -1052.6818
$ws.Range("N2").Value = -1315.4286

$ws.Range("H61").Value = 4335.9707
$ws.Range("I61").Value = 3533.889
$ws.Range("J61").Value = 7429.7144
$ws.Range("K61").Value = 3533.889
$ws.Range("L61").Value = 7429.7144
$ws.Range("M61").Value = -3321.889
$ws.Range("N61").Value = -7853.7144

$ws.Range("H116").Value = 1147.2759
$ws.Range("I116").Value = 1165.6818
$ws.Range("J116").Value = 1089.4286
$ws.Range("K116").Value = 1165.6818
$ws.Range("L116").Value = 1089.4286
$ws.Range("M116").Value = 1128.3182
$ws.Range("N116").Value = -5677.4286

$ws.Range("H136").Value = 4335.9707
$ws.Range("I136").Value = 3533.889
$ws.Range("J136").Value = 7429.7144
$ws.Range("K136").Value = 10601.667
$ws.Range("L136").Value = 22289.1432
$ws.Range("M136").Value = -8051.667000000001
$ws.Range("N136").Value = -27389.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1147.2759
$ws.Range("I3").Value = 1165.6818
$ws.Range("J3").Value = 1089.4286
$ws.Range("K3").Value = 1165.6818
$ws.Range("L3").Value = 1089.4286
$ws.Range("M3").Value = -1051.6818
$ws.Range("N3").Value = -1317.4286

$ws.Range("H94").Value = 1174.5264
$ws.Range("I94").Value = 633.1818
$ws.Range("J94").Value = 1918.875
$ws.Range("K94").Value = 633.1818
$ws.Range("L94").Value = 1918.875
$ws.Range("M94").Value = -182.1818
$ws.Range("N94").Value = -2820.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3557.1562
$ws.Range("I134").Value = 3034.4285
$ws.Range("J134").Value = 3963.7222
$ws.Range("K134").Value = 9103.2855
$ws.Range("L134").Value = 11891.1666
$ws.Range("M134").Value = -6568.2855
$ws.Range("N134").Value = -16961.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23107.924
$ws.Range("I2").Value = 150007.5
$ws.Range("J2").Value = 35.272728
$ws.Range("K2").Value = 900045
$ws.Range("L2").Value = 211.636368
$ws.Range("M2").Value = -899932
$ws.Range("N2").Value = -437.636368

$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9312

$ws.Range("H58").Value = 2942.8572
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2942.8572
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8828.571599999999
$ws.Range("N58").Value = -9084.571599999999

$ws.Range("H64").Value = 47621976
$ws.Range("I64").Value = 200001180
$ws.Range("J64").Value = 3474.6875
$ws.Range("K64").Value = 600003540
$ws.Range("L64").Value = 10424.0625
$ws.Range("M64").Value = -600003270
$ws.Range("N64").Value = -10964.0625

$ws.Range("H67").Value = 47621976
$ws.Range("I67").Value = 200001180
$ws.Range("J67").Value = 3474.6875
$ws.Range("K67").Value = 600003540
$ws.Range("L67").Value = 10424.0625
$ws.Range("M67").Value = -600002604
$ws.Range("N67").Value = -12296.0625

$ws.Range("H92").Value = 716.9091
$ws.Range("I92").Value = 632.5
$ws.Range("J92").Value = 765.1429000000001
$ws.Range("K92").Value = 1897.5
$ws.Range("L92").Value = 2295.4287
$ws.Range("M92").Value = -649.5
$ws.Range("N92").Value = -4791.4287

$ws.Range("H94").Value = 10000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352

$ws.Range("H107").Value = 1373.5
$ws.Range("I107").Value = 377.5
$ws.Range("J107").Value = 1539.5
$ws.Range("K107").Value = 1132.5
$ws.Range("L107").Value = 4618.5
$ws.Range("M107").Value = 787.5
$ws.Range("N107").Value = -8458.5

$ws.Range("H129").Value = 2214.8572
$ws.Range("I129").Value = 2013.3334
$ws.Range("J129").Value = 2366
$ws.Range("K129").Value = 6040.0002
$ws.Range("L129").Value = 7098
$ws.Range("M129").Value = -1040.0002
$ws.Range("N129").Value = -17098

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 11281.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 11281.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 33844.5
$ws.Range("N136").Value = -38944.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1097.7142
$ws.Range("I93").Value = 825
$ws.Range("J93").Value = 1461.3334
$ws.Range("K93").Value = 825
$ws.Range("L93").Value = 1461.3334
$ws.Range("M93").Value = 423
$ws.Range("N93").Value = -3957.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45939.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 45939.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 45939.5
$ws.Range("N46").Value = -46401.5

$ws.Range("H134").Value = 45939.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 45939.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 137818.5
$ws.Range("N134").Value = -142888.5

Write-Host "Applied Pandaemonium_Profits.xlsx updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
